# Apply updated cryptocurrency price/volume data to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = "68.794.60"; E = "  -0.53%  " },
    @{ Row = 3; D = "3.838.21"; E = "  +2.31%  " },
    @{ Row = 4; D = $null; E = "  -0.01%  " },
    @{ Row = 5; D = "601.16"; E = "  -0.19%  " },
    @{ Row = 6; D = "161.98"; E = "  -3.01%  " },
    @{ Row = 7; D = "3.839.20"; E = "  +2.37%  " },
    @{ Row = 8; D = $null; E = "  -0.01%  " },
    @{ Row = 9; D = $null; E = "  -1.51%  " },
    @{ Row = 10; D = "0.168"; E = "  -1.18%  " },
    @{ Row = 11; D = "6.30"; E = "  -1.34%  " },
    @{ Row = 12; D = "0.459"; E = "  -0.03%  " },
    @{ Row = 13; D = "36.84"; E = "  -3.13%  " },
    @{ Row = 14; D = $null; E = "  -2.09%  " },
    @{ Row = 15; D = "4.479.02"; E = "  +2.24%  " },
    @{ Row = 16; D = "3.834.06"; E = "  +2.29%  " },
    @{ Row = 17; D = "68.915.69"; E = "  -0.37%  " },
    @{ Row = 18; D = "7.52"; E = "  +2.24%  " },
    @{ Row = 19; D = $null; E = "  -0.12%  " },
    @{ Row = 20; D = "17.11"; E = "  -1.58%  " },
    @{ Row = 21; D = "11.33"; E = "  +0.01%  " },
    @{ Row = 22; D = "484.95"; E = "  -1.75%  " },
    @{ Row = 23; D = "0.718"; E = "  -1.53%  " },
    @{ Row = 24; D = $null; E = "  +3.02%  " },
    @{ Row = 25; D = "83.94"; E = "  -1.04%  " },
    @{ Row = 26; D = "2.24"; E = "  -2.80%  " },
    @{ Row = 27; D = "12.07"; E = "  -2.07%  " },
    @{ Row = 28; D = $null; E = "  -0.10%  " },
    @{ Row = 29; D = "9.97"; E = "  -1.47%  " },
    @{ Row = 30; D = $null; E = "  -1.04%  " },
    @{ Row = 31; D = "7.92"; E = "  -2.25%  " },
    @{ Row = 32; D = "3.987.08"; E = "  +2.33%  " },
    @{ Row = 33; D = $null; E = "  -4.27%  " },
    @{ Row = 34; D = "32.13"; E = "  +1.66%  " },
    @{ Row = 35; D = "3.783.56"; E = "  +2.56%  " },
    @{ Row = 36; D = "0.107"; E = "  -1.67%  " },
    @{ Row = 37; D = "1.03"; E = "  +0.99%  " },
    @{ Row = 38; D = $null; E = "  +3.15%  " },
    @{ Row = 39; D = "5.90"; E = "  -1.61%  " },
    @{ Row = 40; D = $null; E = "  -0.04%  " },
    @{ Row = 41; D = "0.319"; E = "  -1.95%  " },
    @{ Row = 42; D = "436.13"; E = "  +1.48%  " },
    @{ Row = 43; D = "2.96"; E = "  -1.44%  " },
    @{ Row = 44; D = "48.49"; E = "  -0.70%  " },
    @{ Row = 45; D = $null; E = "  -0.91%  " },
    @{ Row = 47; D = "8.38"; E = "  -1.23%  " },
    @{ Row = 48; D = "143.48"; E = "  +1.69%  " },
    @{ Row = 49; D = "2.829.90"; E = "  +1.15%  " },
    @{ Row = 50; D = $null; E = "  +1.98%  " },
    @{ Row = 51; D = "25.77"; E = "  +10.85%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        $dCell = $ws.Cells.Item($u.Row, 4)
        # The Price column holds numeric-looking text (e.g. "601.16" or
        # "3.838.21" with thousands separators baked in as literal dots).
        # Writing that straight to .Value would let the COM layer coerce it
        # to a real number (and mangle multi-dot values / introduce float
        # rounding), so force text formatting first, then strip the format
        # override back off so the cell's style index is left untouched,
        # matching the original (unstyled) data cells.
        $dCell.NumberFormat = "@"
        $dCell.Value = $u.D
        $dCell.ClearFormats()
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}
